$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3156.6924
$ws.Range("I43").Value = 2257
$ws.Range("K43").Value = 2257
$ws.Range("M43").Value = -2188
$ws.Range("H68").Value = 59999
$ws.Range("J68").Value = 59999
$ws.Range("L68").Value = 59999
$ws.Range("N68").Value = -61497
$ws.Range("H71").Value = 59999
$ws.Range("J71").Value = 59999
$ws.Range("L71").Value = 179997
$ws.Range("N71").Value = -187485
$ws.Range("H86").Value = 2575.2144
$ws.Range("I86").Value = 2105
$ws.Range("J86").Value = 3421.6
$ws.Range("K86").Value = 2105
$ws.Range("L86").Value = 3421.6
$ws.Range("M86").Value = -982
$ws.Range("N86").Value = -5667.6
$ws.Range("H89").Value = 2575.2144
$ws.Range("I89").Value = 2105
$ws.Range("J89").Value = 3421.6
$ws.Range("K89").Value = 10525
$ws.Range("L89").Value = 17108
$ws.Range("M89").Value = -4909
$ws.Range("N89").Value = -28340
$ws.Range("H116").Value = 27031
$ws.Range("I116").Value = 29264.092
$ws.Range("J116").Value = 14749
$ws.Range("K116").Value = 29264.092
$ws.Range("L116").Value = 14749
$ws.Range("M116").Value = -25822.092
$ws.Range("N116").Value = -21633

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 8183.6665
$ws.Range("J46").Value = 7526
$ws.Range("L46").Value = 7526
$ws.Range("N46").Value = -8164
$ws.Range("H97").Value = 1605.9048
$ws.Range("I97").Value = 1217
$ws.Range("J97").Value = 2383.7144
$ws.Range("K97").Value = 1217
$ws.Range("L97").Value = 2383.7144
$ws.Range("M97").Value = -721
$ws.Range("N97").Value = -3375.7144
$ws.Range("H132").Value = 1067.6666
$ws.Range("I132").Value = 1067.6666
$ws.Range("K132").Value = 3202.9998
$ws.Range("M132").Value = -672.9998000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 2471
$ws.Range("I26").Value = 2471
$ws.Range("K26").Value = 2471
$ws.Range("M26").Value = -2179
$ws.Range("H96").Value = 9850
$ws.Range("I96").Value = 9850
$ws.Range("K96").Value = 9850
$ws.Range("M96").Value = -7104
$ws.Range("H99").Value = 1166
$ws.Range("I99").Value = 1166
$ws.Range("K99").Value = 1166
$ws.Range("M99").Value = 332
$ws.Range("H105").Value = 1682.7273
$ws.Range("I105").Value = 1152.6818
$ws.Range("J105").Value = 2742.818
$ws.Range("K105").Value = 1152.6818
$ws.Range("L105").Value = 2742.818
$ws.Range("M105").Value = 594.3181999999999
$ws.Range("N105").Value = -6236.818

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 614
$ws.Range("I19").Value = 614
$ws.Range("K19").Value = 614
$ws.Range("M19").Value = -444
$ws.Range("H24").Value = 614
$ws.Range("I24").Value = 614
$ws.Range("K24").Value = 614
$ws.Range("M24").Value = -444
$ws.Range("H94").Value = 3125.7334
$ws.Range("I94").Value = 1892.6666
$ws.Range("J94").Value = 3947.7778
$ws.Range("K94").Value = 1892.6666
$ws.Range("L94").Value = 3947.7778
$ws.Range("M94").Value = -1441.6666
$ws.Range("N94").Value = -4849.7778
$ws.Range("H99").Value = 5655.35
$ws.Range("I99").Value = 4858.5
$ws.Range("J99").Value = 6452.2
$ws.Range("K99").Value = 4858.5
$ws.Range("L99").Value = 6452.2
$ws.Range("M99").Value = -3360.5
$ws.Range("N99").Value = -9448.200000000001
$ws.Range("H105").Value = 33454.082
$ws.Range("I105").Value = 48917.125
$ws.Range("K105").Value = 48917.125
$ws.Range("M105").Value = -47170.125
$ws.Range("H126").Value = 5655.35
$ws.Range("I126").Value = 4858.5
$ws.Range("J126").Value = 6452.2
$ws.Range("K126").Value = 14575.5
$ws.Range("L126").Value = 19356.6
$ws.Range("M126").Value = -12105.5
$ws.Range("N126").Value = -24296.6
$ws.Range("H132").Value = 64614.75
$ws.Range("J132").Value = 3690.2856
$ws.Range("L132").Value = 11070.8568
$ws.Range("N132").Value = -16130.8568
$ws.Range("H134").Value = 3456.5483
$ws.Range("I134").Value = 2933.3914
$ws.Range("K134").Value = 8800.174199999999
$ws.Range("M134").Value = -6265.174199999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 319.08334
$ws.Range("J23").Value = 326
$ws.Range("L23").Value = 978
$ws.Range("N23").Value = -1448
$ws.Range("H46").Value = 1750
$ws.Range("I46").Value = 1440
$ws.Range("J46").Value = 2525
$ws.Range("K46").Value = 4320
$ws.Range("L46").Value = 7575
$ws.Range("M46").Value = -4229
$ws.Range("N46").Value = -7757
$ws.Range("H60").Value = 2041.5
$ws.Range("I60").Value = 769.1667
$ws.Range("J60").Value = 3950
$ws.Range("K60").Value = 2307.5001
$ws.Range("L60").Value = 11850
$ws.Range("M60").Value = -2056.5001
$ws.Range("N60").Value = -12352
$ws.Range("H87").Value = 16796.4
$ws.Range("I87").Value = 8745.5
$ws.Range("K87").Value = 26236.5
$ws.Range("M87").Value = -24988.5
$ws.Range("H90").Value = 16796.4
$ws.Range("I90").Value = 8745.5
$ws.Range("K90").Value = 78709.5
$ws.Range("M90").Value = -72469.5
$ws.Range("H107").Value = 2322.1
$ws.Range("J107").Value = 1149.1666
$ws.Range("L107").Value = 3447.4998
$ws.Range("N107").Value = -7287.4998
$ws.Range("H131").Value = 108740.4
$ws.Range("J131").Value = 1921.4517
$ws.Range("L131").Value = 5764.355100000001
$ws.Range("N131").Value = -15844.3551

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 40000
$ws.Range("J47").Value = 40000
$ws.Range("L47").Value = 40000
$ws.Range("N47").Value = -41136
$ws.Range("H70").Value = 4460.8237
$ws.Range("I70").Value = 4367.2
$ws.Range("J70").Value = 4594.5713
$ws.Range("K70").Value = 4367.2
$ws.Range("L70").Value = 4594.5713
$ws.Range("M70").Value = -4097.2
$ws.Range("N70").Value = -5134.5713
$ws.Range("H73").Value = 4460.8237
$ws.Range("I73").Value = 4367.2
$ws.Range("J73").Value = 4594.5713
$ws.Range("K73").Value = 4367.2
$ws.Range("L73").Value = 4594.5713
$ws.Range("M73").Value = -3431.2
$ws.Range("N73").Value = -6466.5713
$ws.Range("H132").Value = 2694.6667
$ws.Range("I132").Value = 2354.1538
$ws.Range("J132").Value = 3580
$ws.Range("K132").Value = 7062.4614
$ws.Range("L132").Value = 10740
$ws.Range("M132").Value = -4532.4614
$ws.Range("N132").Value = -15800

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5067.143
$ws.Range("I46").Value = 1780.2
$ws.Range("J46").Value = 13284.5
$ws.Range("K46").Value = 1780.2
$ws.Range("L46").Value = 13284.5
$ws.Range("M46").Value = -1592.2
$ws.Range("N46").Value = -13660.5
$ws.Range("H68").Value = 4092.611
$ws.Range("I68").Value = 3379.182
$ws.Range("J68").Value = 5213.7144
$ws.Range("K68").Value = 3379.182
$ws.Range("L68").Value = 5213.7144
$ws.Range("M68").Value = -2630.182
$ws.Range("N68").Value = -6711.7144
$ws.Range("H71").Value = 4092.611
$ws.Range("I71").Value = 3379.182
$ws.Range("J71").Value = 5213.7144
$ws.Range("K71").Value = 16895.91
$ws.Range("L71").Value = 26068.572
$ws.Range("M71").Value = -13151.91
$ws.Range("N71").Value = -33556.572
$ws.Range("H82").Value = 1428.0667
$ws.Range("I82").Value = 1577.2
$ws.Range("J82").Value = 1129.8
$ws.Range("K82").Value = 1577.2
$ws.Range("L82").Value = 1129.8
$ws.Range("M82").Value = -1216.2
$ws.Range("N82").Value = -1851.8
$ws.Range("H85").Value = 1428.0667
$ws.Range("I85").Value = 1577.2
$ws.Range("J85").Value = 1129.8
$ws.Range("K85").Value = 1577.2
$ws.Range("L85").Value = 1129.8
$ws.Range("M85").Value = -329.2
$ws.Range("N85").Value = -3625.8
$ws.Range("H122").Value = 2886.3333
$ws.Range("I122").Value = 2706
$ws.Range("K122").Value = 8118
$ws.Range("M122").Value = -5668
$ws.Range("H128").Value = 149999
$ws.Range("J128").Value = 149999
$ws.Range("L128").Value = 149999
$ws.Range("N128").Value = -159959
$ws.Range("H132").Value = 4922.8
$ws.Range("I132").Value = 4922.8
$ws.Range("K132").Value = 14768.4
$ws.Range("M132").Value = -12238.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 31000
$ws.Range("J119").Value = 31000
$ws.Range("L119").Value = 31000
$ws.Range("N119").Value = -40676
